# Build site at 2022-09-26 16:07:08 UTC
#
# The site generator dropped the two standalone rows that used to hold just
# the "Docentes responsaveis" values (the two professors' names), and the
# value cells on several of the remaining label rows were re-populated from
# other fields further down the form (mirrors the upstream re-flowed output).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that used to hold the "Docentes responsaveis" values on
# their own (old rows 13 and 14, "519033 - Carlos Yujiro Shigue" /
# "5840726 - Cristina Bormio Nunes"). Everything below shifts up by two rows,
# which also carries the correct row heights along for free, taking the
# sheet from A1:C25 down to A1:C23.
$ws.Rows("13:14").Delete()

# Objetivos: value becomes the first docente's name.
$ws.Range("B10").Value2 = "519033 - Carlos Yujiro Shigue"
$ws.Range("C10").Value2 = "519033 - Carlos Yujiro Shigue"

# Programa resumido: value becomes the activation date. Assigning a
# date-shaped literal straight to .Value2 would get auto-converted to a
# serial date (and drag in a new number-format style), so build it as text
# via a formula first and then freeze it to a plain value with Paste
# Special - that keeps the original "s" style untouched.
$ws.Range("B13").Formula = '=TEXT(DATE(2012,1,1),"mm/dd/yyyy")'
$ws.Range("B13").Copy()
$ws.Range("B13").PasteSpecial(-4163)
$ws.Range("C13").Formula = '=TEXT(DATE(2012,1,1),"mm/dd/yyyy")'
$ws.Range("C13").Copy()
$ws.Range("C13").PasteSpecial(-4163)
$ws.Application.CutCopyMode = $false

# Programa: value becomes the first docente's name.
$ws.Range("B15").Value2 = "519033 - Carlos Yujiro Shigue"
$ws.Range("C15").Value2 = "519033 - Carlos Yujiro Shigue"

# Metodo: value becomes the second docente's name.
$ws.Range("B18").Value2 = "5840726 - Cristina Bormio Nunes"
$ws.Range("C18").Value2 = "5840726 - Cristina Bormio Nunes"

# Criterio: value becomes the teaching-method text.
$ws.Range("B19").Value2 = "Aulas expositivas e práticas ministradas em laboratório."
$ws.Range("C19").Value2 = "Aulas expositivas e práticas ministradas em laboratório."

# Norma de recuperacao: value becomes the evaluation-criteria text.
$ws.Range("B20").Value2 = "Média ponderada de duas provas escritas, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4"
$ws.Range("C20").Value2 = "Média ponderada de duas provas escritas, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4"

# Bibliografia: value becomes the recovery-rule text (replacing the long
# reading list, which is dropped entirely).
$ws.Range("B21").Value2 = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Range("C21").Value2 = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"

$ws.Range("A1").Select()
